$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell values (content corrections) ---
$ws.Range("B10").Value = 'Apresentar uma visão geral da química dos elementos e de seus compostos enfatizando as correlações entre as propriedades físicas e químicas com os aspectos estruturais e de ligação, os métodos de obtenção em laboratório e indústria, além das principais propriedades e aplicações.'
$ws.Range("C10").Value = 'Apresentar uma visão geral da química dos elementos e de seus compostos enfatizando as correlações entre as propriedades físicas e químicas com os aspectos estruturais e de ligação, os métodos de obtenção em laboratório e indústria, além das principais propriedades e aplicações.'
$ws.Range("B13").Value = '5840963 - Daniela Camargo Vernilli'
$ws.Range("C13").Value = '5840963 - Daniela Camargo Vernilli'
$ws.Range("A14").Value = 'Programa resumido:'
$ws.Range("B14").Value = 'Ocorrência, obtenção, estrutura, propriedades e aplicações de elementos metálicos e não-metálicos; moléculas poliatômicas; compostos halogenados e das famílias do oxigênio, nitrogênio, carbono e boro; compostos oxigenados. Processos industriais de fabricação.'
$ws.Range("C14").Value = 'Ocorrência, obtenção, estrutura, propriedades e aplicações de elementos metálicos e não-metálicos; moléculas poliatômicas; compostos halogenados e das famílias do oxigênio, nitrogênio, carbono e boro; compostos oxigenados. Processos industriais de fabricação.'
$ws.Range("A15").Value = 'Short syllabus:'
$ws.Range("B15").Value = 'Occurrence, obtaining, structure, properties and applications of metallic and non-metallic elements; polyatomic molecules; halogenated compounds and the oxygen, nitrogen, carbon and boron families; oxygenated compounds. Industrial manufacturing processes.'
$ws.Range("C15").Value = 'Occurrence, obtaining, structure, properties and applications of metallic and non-metallic elements; polyatomic molecules; halogenated compounds and the oxygen, nitrogen, carbon and boron families; oxygenated compounds. Industrial manufacturing processes.'
$ws.Range("A16").Value = 'Programa:'
$ws.Range("B16").Value = 'Ocorrência, obtenção, estrutura, propriedades e aplicações de elementos não-metálicos: gases nobres, hidrogênio molecular, halogênios, oxigênio molecular, ozônio e nitrogênio molecular; semimetais; metais alcalinos, alcalinos-terrosos e de transição; moléculas poliatômicas e espécies catenadas de: enxofre, fósforo e carbono; compostos halogenados e das famílias do oxigênio, nitrogênio, carbono e boro; compostos oxigenados: óxidos e oxicompostos. Processos industriais de fabricação dos principais insumos químicos e materiais.'
$ws.Range("C16").Value = 'Ocorrência, obtenção, estrutura, propriedades e aplicações de elementos não-metálicos: gases nobres, hidrogênio molecular, halogênios, oxigênio molecular, ozônio e nitrogênio molecular; semimetais; metais alcalinos, alcalinos-terrosos e de transição; moléculas poliatômicas e espécies catenadas de: enxofre, fósforo e carbono; compostos halogenados e das famílias do oxigênio, nitrogênio, carbono e boro; compostos oxigenados: óxidos e oxicompostos. Processos industriais de fabricação dos principais insumos químicos e materiais.'
$ws.Range("A17").Value = 'Syllabus:'
$ws.Range("A18").Value = 'Avaliação:'
$ws.Range("A19").Value = 'Método:'
$ws.Range("A20").Value = 'Critério:'
$ws.Range("A21").Value = 'Norma de recuperação:'
$ws.Range("A22").Value = 'Bibliografia:'

# --- New cells: copy formatting from an existing same-column cell, then set value ---
$ws.Range("B3").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("B17").Value = 'Occurrence, obtaining, structure, properties and applications of non-metallic elements: noble gases, molecular hydrogen, halogens, molecular oxygen, ozone and molecular nitrogen; semimetals; alkali, alkaline earth and transition metals; polyatomic molecules and catenated species of: sulfur, phosphorus and carbon; halogenated compounds and the oxygen, nitrogen, carbon and boron families; oxygenated compounds: oxides and oxycompounds. Industrial manufacturing processes of the main chemical inputs and materials.'
$ws.Range("C3").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C17").Value = 'Occurrence, obtaining, structure, properties and applications of non-metallic elements: noble gases, molecular hydrogen, halogens, molecular oxygen, ozone and molecular nitrogen; semimetals; alkali, alkaline earth and transition metals; polyatomic molecules and catenated species of: sulfur, phosphorus and carbon; halogenated compounds and the oxygen, nitrogen, carbon and boron families; oxygenated compounds: oxides and oxycompounds. Industrial manufacturing processes of the main chemical inputs and materials.'
$ws.Range("B3").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("B22").Value = 'QUAGLIANO, J. V.; VALLARINO, L. Química, Guanabara Koogan, 1973.
LEE, J. D. Química Inorgânica, Editora Edgard Blücher, 1999.
GREENWOOD, N. N.; EARNSHAW, A. Chemistry of the Elements, Butterworth Heinemann, 1997.
SHRIVER, D. F.; ATKINS, P. W.; LANGFORD, G. H. Inorganic Chemistry, Oxford University Press, 1994.
PORTERFIELD, W. W. Inorganic Chemistry: a Unified approach, Addison Wesley Heading, 1984.
BUCHNER, W.; SCHLIEBS, R.; WINTER, G.; BUCHEL, K. H. Industrial Inorganic Chemistry, VCH, 1989.'
$ws.Range("C3").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = 'QUAGLIANO, J. V.; VALLARINO, L. Química, Guanabara Koogan, 1973.
LEE, J. D. Química Inorgânica, Editora Edgard Blücher, 1999.
GREENWOOD, N. N.; EARNSHAW, A. Chemistry of the Elements, Butterworth Heinemann, 1997.
SHRIVER, D. F.; ATKINS, P. W.; LANGFORD, G. H. Inorganic Chemistry, Oxford University Press, 1994.
PORTERFIELD, W. W. Inorganic Chemistry: a Unified approach, Addison Wesley Heading, 1984.
BUCHNER, W.; SCHLIEBS, R.; WINTER, G.; BUCHEL, K. H. Industrial Inorganic Chemistry, VCH, 1989.'
$ws.Range("A3").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("A23").Value = 'Requisitos:'
$ws.Range("B3").Copy()
$ws.Range("B24").PasteSpecial(-4122)
$ws.Range("B24").Value = 'LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito fraco)
'
$ws.Range("C3").Copy()
$ws.Range("C24").PasteSpecial(-4122)
$ws.Range("C24").Value = 'LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito fraco)
'

# --- Clear cells that no longer hold content ---
$ws.Range("A13").Clear()
$ws.Range("B18").Clear()
$ws.Range("C18").Clear()
$ws.Range("B23").Clear()
$ws.Range("C23").Clear()

# --- Adjust row heights to match final layout ---
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).AutoFit()
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 120
$ws.Rows.Item(23).AutoFit()
$ws.Rows.Item(24).RowHeight = 30

$excel.CutCopyMode = 0
